$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new data row (row 6) for "PEDOMAN UMUM TATA NASKAH DINAS", mirroring the
# existing dropdown entries (copy formatting from row 2, then overwrite values).
$ws.Range("A2:J2").Copy($ws.Range("A6:J6"))
$ws.Rows(6).RowHeight = 68.5

$ws.Range("A6").Value = "peraturan_menteri_lembaga"
$ws.Range("C6").Value = "5 TAHUN 2021"
$ws.Range("D6").Value = "PERATURAN ARSIP NASIONAL REPUBLIK INDONESIA"
$ws.Range("B6").Value = "PEDOMAN UMUM TATA NASKAH DINAS"
$ws.Range("I6").Value = "BERITA NEGARA REPUBLIK INDONESIA TAHUN 2021 NOMOR 758"
$ws.Range("E6").Value = "peraturan lembaga"
$ws.Range("F6").Value = "jakarta"
$ws.Range("G6").Value = 44372
$ws.Range("H6").Value = 44372
$ws.Range("J6").Value = "Umum"

[void]$ws.Range("D6").Select()
